# Aggiornamento dati: aggiunta righe per le date 15, 16, 17 marzo
# (date seriali Excel 44301, 44302, 44303) in coda alla tabella.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Nuovi dati da accodare: riga, data (seriale), nuovi positivi, somma mobile 7gg, somma mobile per 100mila abitanti
$newRows = @(
    @(227, 44301, 8,  60, 234.842850992211),
    @(228, 44302, 10, 55, 215.2726134095268),
    @(229, 44303, 11, 51, 199.6164233433794)
)

foreach ($entry in $newRows) {
    $r        = $entry[0]
    $dateVal  = $entry[1]
    $bVal     = $entry[2]
    $cVal     = $entry[3]
    $dVal     = $entry[4]

    # Copia lo stile della colonna A dalla riga precedente (contiene lo stile data)
    # cosi' la nuova cella mantiene la stessa formattazione, poi si imposta il valore.
    $srcA = $ws.Cells.Item($r - 1, 1)
    $dstA = $ws.Cells.Item($r, 1)
    $srcA.Copy($dstA)
    $dstA.Value2 = $dateVal

    $ws.Cells.Item($r, 2).Value2 = $bVal
    $ws.Cells.Item($r, 3).Value2 = $cVal
    $ws.Cells.Item($r, 4).Value2 = $dVal
}
